$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.5386712677833714
$ws.Cells.Item(2, 3).Value = 0.1051590893890264
$ws.Cells.Item(2, 4).Value = 0.04124891684684684
$ws.Cells.Item(2, 5).Value = 0.09797392431061525
$ws.Cells.Item(2, 6).Value = 0.8658998088534631
$ws.Cells.Item(2, 8).Value = 0.07973214163530429
$ws.Cells.Item(2, 9).Value = 0.7959050165789705
$ws.Cells.Item(2, 11).Value = 0.3331153853719968
$ws.Cells.Item(2, 12).Value = 0.2056283115882991
$ws.Cells.Item(2, 14).Value = 1.617569122752696
$ws.Cells.Item(2, 15).Value = 3.05839678202031

$ws.Cells.Item(3, 2).Value = 0.4977341740913062
$ws.Cells.Item(3, 3).Value = 0.1035833427660435
$ws.Cells.Item(3, 4).Value = 0.03904359781390809
$ws.Cells.Item(3, 5).Value = 0.09743708183724209
$ws.Cells.Item(3, 6).Value = 0.8657033429486845
$ws.Cells.Item(3, 8).Value = 0.07973214163530429
$ws.Cells.Item(3, 9).Value = 0.8011007215043549
$ws.Cells.Item(3, 11).Value = 0.2959267575989486
$ws.Cells.Item(3, 12).Value = 0.1982594525952663
$ws.Cells.Item(3, 14).Value = 1.633924489781511
$ws.Cells.Item(3, 15).Value = 3.070715944059202

$ws.Cells.Item(4, 2).Value = 0.4727407509763282
$ws.Cells.Item(4, 3).Value = 0.1026049566044236
$ws.Cells.Item(4, 4).Value = 0.03767506608364357
$ws.Cells.Item(4, 5).Value = 0.09715595846729386
$ws.Cells.Item(4, 6).Value = 0.86606032163413
$ws.Cells.Item(4, 8).Value = 0.07973214163530429
$ws.Cells.Item(4, 9).Value = 0.8046630018676808
$ws.Cells.Item(4, 11).Value = 0.2731019816422133
$ws.Cells.Item(4, 12).Value = 0.1938435960835108
$ws.Cells.Item(4, 14).Value = 1.644483976683105
$ws.Cells.Item(4, 15).Value = 3.079838884716651

$ws.Cells.Item(5, 2).Value = 0.4625921552019747
$ws.Cells.Item(5, 3).Value = 0.1022035415682012
$ws.Cells.Item(5, 4).Value = 0.03711376696978874
$ws.Cells.Item(5, 5).Value = 0.09705361697474757
$ws.Cells.Item(5, 6).Value = 0.8663259775372794
$ws.Cells.Item(5, 8).Value = 0.07973214163530429
$ws.Cells.Item(5, 9).Value = 0.8062082719799015
$ws.Cells.Item(5, 11).Value = 0.2638035503759539
$ws.Cells.Item(5, 12).Value = 0.1920715128021158
$ws.Cells.Item(5, 14).Value = 1.648917129387678
$ws.Cells.Item(5, 15).Value = 3.083948746751872

$ws.Cells.Item(6, 2).Value = 0.4609092090261413
$ws.Cells.Item(6, 3).Value = 0.1021367235026247
$ws.Cells.Item(6, 4).Value = 0.0370203462697134
$ws.Cells.Item(6, 5).Value = 0.09703736188987477
$ws.Cells.Item(6, 6).Value = 0.8663773508395991
$ws.Cells.Item(6, 8).Value = 0.07973214163530429
$ws.Cells.Item(6, 9).Value = 0.8064705181707481
$ws.Cells.Item(6, 11).Value = 0.2622597409564378
$ws.Cells.Item(6, 12).Value = 0.1917789181870972
$ws.Cells.Item(6, 14).Value = 1.649661107784145
$ws.Cells.Item(6, 15).Value = 3.084654878968749

$ws.Cells.Item(7, 2).Value = 0.4726037350574757
$ws.Cells.Item(7, 3).Value = 0.1025995539479183
$ws.Cells.Item(7, 4).Value = 0.03766751079787412
$ws.Cells.Item(7, 5).Value = 0.09715452875114394
$ws.Cells.Item(7, 6).Value = 0.8660634176283253
$ws.Cells.Item(7, 8).Value = 0.07973214163530429
$ws.Cells.Item(7, 9).Value = 0.8046834628297326
$ws.Cells.Item(7, 11).Value = 0.2729765674920372
$ws.Cells.Item(7, 12).Value = 0.1938195860272032
$ws.Cells.Item(7, 14).Value = 1.644543237036874
$ws.Cells.Item(7, 15).Value = 3.079892723548085

$ws.Cells.Item(8, 2).Value = 0.5245270623259159
$ws.Cells.Item(8, 3).Value = 0.1046180392472067
$ws.Cells.Item(8, 4).Value = 0.04049153719944343
$ws.Cells.Item(8, 5).Value = 0.09777876813336306
$ws.Cells.Item(8, 6).Value = 0.8657329883881104
$ws.Cells.Item(8, 8).Value = 0.07973214163530429
$ws.Cells.Item(8, 9).Value = 0.7976192937273936
$ws.Cells.Item(8, 11).Value = 0.3202911507577824
$ws.Cells.Item(8, 12).Value = 0.2030650309881565
$ws.Cells.Item(8, 14).Value = 1.623101113967246
$ws.Cells.Item(8, 15).Value = 3.062320980859909

$ws.Cells.Item(9, 2).Value = 0.6274510899455379
$ws.Cells.Item(9, 3).Value = 0.1084893561268885
$ws.Cells.Item(9, 4).Value = 0.04591390937511619
$ws.Cells.Item(9, 5).Value = 0.09938689841987625
$ws.Cells.Item(9, 6).Value = 0.8688720300209809
$ws.Cells.Item(9, 8).Value = 0.07973214163530429
$ws.Cells.Item(9, 9).Value = 0.7867171728087605
$ws.Cells.Item(9, 11).Value = 0.4131286422376093
$ws.Cells.Item(9, 12).Value = 0.2220544979290793
$ws.Cells.Item(9, 14).Value = 1.585155804068846
$ws.Cells.Item(9, 15).Value = 3.040226325272158

$ws.Cells.Item(10, 2).Value = 0.7037154811095263
$ws.Cells.Item(10, 3).Value = 0.1112799672338554
$ws.Cells.Item(10, 4).Value = 0.04982651940346727
$ws.Cells.Item(10, 5).Value = 0.1008016338701339
$ws.Cells.Item(10, 6).Value = 0.873485142868283
$ws.Cells.Item(10, 8).Value = 0.07973214163530429
$ws.Cells.Item(10, 9).Value = 0.7805042750362077
$ws.Cells.Item(10, 11).Value = 0.4813499111045871
$ws.Cells.Item(10, 12).Value = 0.2365275647395748
$ws.Cells.Item(10, 14).Value = 1.559775730341777
$ws.Cells.Item(10, 15).Value = 3.031525887103271

$ws.Cells.Item(11, 2).Value = 0.7385453508412354
$ws.Cells.Item(11, 3).Value = 0.1125377160015546
$ws.Cells.Item(11, 4).Value = 0.05159085897007287
$ws.Cells.Item(11, 5).Value = 0.1014957141737547
$ws.Cells.Item(11, 6).Value = 0.8760844041170586
$ws.Cells.Item(11, 8).Value = 0.07973214163530429
$ws.Cells.Item(11, 9).Value = 0.7780676722789011
$ws.Cells.Item(11, 11).Value = 0.5123847872752663
$ws.Cells.Item(11, 12).Value = 0.2432245889232547
$ws.Cells.Item(11, 14).Value = 1.548771592419325
$ws.Cells.Item(11, 15).Value = 3.029202715278188

$ws.Cells.Item(12, 2).Value = 0.7517535861377951
$ws.Cells.Item(12, 3).Value = 0.1130122920420931
$ws.Cells.Item(12, 4).Value = 0.05225671512011587
$ws.Cells.Item(12, 5).Value = 0.1017657882382466
$ws.Cells.Item(12, 6).Value = 0.8771406247382529
$ws.Cells.Item(12, 8).Value = 0.07973214163530429
$ws.Cells.Item(12, 9).Value = 0.7772009957414667
$ws.Cells.Item(12, 11).Value = 0.5241365198425854
$ws.Cells.Item(12, 12).Value = 0.2457767771812343
$ws.Cells.Item(12, 14).Value = 1.544682453440031
$ws.Cells.Item(12, 15).Value = 3.028557945559186

$ws.Cells.Item(13, 2).Value = 0.7489081269085887
$ws.Cells.Item(13, 3).Value = 0.1129101598233362
$ws.Cells.Item(13, 4).Value = 0.05211341209227527
$ws.Cells.Item(13, 5).Value = 0.1017073012718761
$ws.Cells.Item(13, 6).Value = 0.8769099504141309
$ws.Cells.Item(13, 8).Value = 0.07973214163530429
$ws.Cells.Item(13, 9).Value = 0.7773851592540097
$ws.Cells.Item(13, 11).Value = 0.5216056058191896
$ws.Cells.Item(13, 12).Value = 0.2452264001988311
$ws.Cells.Item(13, 14).Value = 1.545559657896682
$ws.Cells.Item(13, 15).Value = 3.028686359774269

$ws.Cells.Item(14, 2).Value = 0.7396316254565534
$ws.Cells.Item(14, 3).Value = 0.1125767939738012
$ws.Cells.Item(14, 4).Value = 0.05164568483122878
$ws.Cells.Item(14, 5).Value = 0.1015177883393719
$ws.Cells.Item(14, 6).Value = 0.8761698586956754
$ws.Cells.Item(14, 8).Value = 0.07973214163530429
$ws.Cells.Item(14, 9).Value = 0.7779952478446859
$ws.Cells.Item(14, 11).Value = 0.5133516228709993
$ws.Cells.Item(14, 12).Value = 0.2434342354853953
$ws.Cells.Item(14, 14).Value = 1.548433615247713
$ws.Cells.Item(14, 15).Value = 3.029144961230287

$ws.Cells.Item(15, 2).Value = 0.7339519419436442
$ws.Cells.Item(15, 3).Value = 0.1123723749368679
$ws.Cells.Item(15, 4).Value = 0.05135889309566721
$ws.Cells.Item(15, 5).Value = 0.1014026485114137
$ws.Cells.Item(15, 6).Value = 0.87572589740806
$ws.Cells.Item(15, 8).Value = 0.07973214163530429
$ws.Cells.Item(15, 9).Value = 0.7783762384658885
$ws.Cells.Item(15, 11).Value = 0.5082957368846337
$ws.Cells.Item(15, 12).Value = 0.2423385857288167
$ws.Cells.Item(15, 14).Value = 1.550204141776784
$ws.Cells.Item(15, 15).Value = 3.029456464170408

$ws.Cells.Item(16, 2).Value = 0.7014419558196607
$ws.Cells.Item(16, 3).Value = 0.1111975328245194
$ws.Cells.Item(16, 4).Value = 0.04971090081737373
$ws.Cells.Item(16, 5).Value = 0.1007572885719448
$ws.Cells.Item(16, 6).Value = 0.873325345270672
$ws.Cells.Item(16, 8).Value = 0.07973214163530429
$ws.Cells.Item(16, 9).Value = 0.7806713520751742
$ws.Cells.Item(16, 11).Value = 0.4793216714207915
$ws.Cells.Item(16, 12).Value = 0.2360921667980733
$ws.Cells.Item(16, 14).Value = 1.560505773939872
$ws.Cells.Item(16, 15).Value = 3.031710597432664

$ws.Cells.Item(17, 2).Value = 0.681532637566562
$ws.Cells.Item(17, 3).Value = 0.1104737895485854
$ws.Cells.Item(17, 4).Value = 0.04869591416968433
$ws.Cells.Item(17, 5).Value = 0.1003743021107226
$ws.Cells.Item(17, 6).Value = 0.8719808859547271
$ws.Cells.Item(17, 8).Value = 0.07973214163530429
$ws.Cells.Item(17, 9).Value = 0.7821791189718752
$ws.Cells.Item(17, 11).Value = 0.4615467685254089
$ws.Cells.Item(17, 12).Value = 0.2322891009502683
$ws.Cells.Item(17, 14).Value = 1.566964180069243
$ws.Cells.Item(17, 15).Value = 3.033512048426388

$ws.Cells.Item(18, 2).Value = 0.6700942471656504
$ws.Cells.Item(18, 3).Value = 0.1100564102531081
$ws.Cells.Item(18, 4).Value = 0.04811066097736472
$ws.Cells.Item(18, 5).Value = 0.1001587740178742
$ws.Cells.Item(18, 6).Value = 0.8712547240074144
$ws.Cells.Item(18, 8).Value = 0.07973214163530429
$ws.Cells.Item(18, 9).Value = 0.7830830260889599
$ws.Cells.Item(18, 11).Value = 0.4513232142533354
$ws.Cells.Item(18, 12).Value = 0.2301123346781395
$ws.Cells.Item(18, 14).Value = 1.570729839693519
$ws.Cells.Item(18, 15).Value = 3.034702067418607

$ws.Cells.Item(19, 2).Value = 0.6662236472069765
$ws.Cells.Item(19, 3).Value = 0.1099149043042829
$ws.Cells.Item(19, 4).Value = 0.04791225445403313
$ws.Cells.Item(19, 5).Value = 0.1000866173657435
$ws.Cells.Item(19, 6).Value = 0.8710169556729284
$ws.Cells.Item(19, 8).Value = 0.07973214163530429
$ws.Cells.Item(19, 9).Value = 0.7833953739392854
$ws.Cells.Item(19, 11).Value = 0.4478617298171343
$ws.Cells.Item(19, 12).Value = 0.2293771523636821
$ws.Cells.Item(19, 14).Value = 1.57201357781163
$ws.Cells.Item(19, 15).Value = 3.035131418195959

$ws.Cells.Item(20, 2).Value = 0.6836506864755734
$ws.Cells.Item(20, 3).Value = 0.1105509473930084
$ws.Cells.Item(20, 4).Value = 0.04880411251115646
$ws.Cells.Item(20, 5).Value = 0.1004145796493674
$ws.Cells.Item(20, 6).Value = 0.8721191280870002
$ws.Cells.Item(20, 8).Value = 0.07973214163530429
$ws.Cells.Item(20, 9).Value = 0.7820148185883014
$ws.Cells.Item(20, 11).Value = 0.4634389312974179
$ws.Cells.Item(20, 12).Value = 0.2326928416107279
$ws.Cells.Item(20, 14).Value = 1.566271397761099
$ws.Cells.Item(20, 15).Value = 3.033304356455943

$ws.Cells.Item(21, 2).Value = 0.7423558513435466
$ws.Cells.Item(21, 3).Value = 0.1126747580807574
$ws.Cells.Item(21, 4).Value = 0.05178312925306727
$ws.Cells.Item(21, 5).Value = 0.1015732565439791
$ws.Cells.Item(21, 6).Value = 0.8763852897762732
$ws.Cells.Item(21, 8).Value = 0.07973214163530429
$ws.Cells.Item(21, 9).Value = 0.77781453006817
$ws.Cells.Item(21, 11).Value = 0.5157760361310011
$ws.Cells.Item(21, 12).Value = 0.2439601996477734
$ws.Cells.Item(21, 14).Value = 1.547587350096141
$ws.Cells.Item(21, 15).Value = 3.029003882809207

$ws.Cells.Item(22, 2).Value = 0.7808329971527996
$ws.Cells.Item(22, 3).Value = 0.1140528365424345
$ws.Cells.Item(22, 4).Value = 0.05371689714729655
$ws.Cells.Item(22, 5).Value = 0.1023727140170472
$ws.Cells.Item(22, 6).Value = 0.8795927435220392
$ws.Cells.Item(22, 8).Value = 0.07973214163530429
$ws.Cells.Item(22, 9).Value = 0.7753958691518648
$ws.Cells.Item(22, 11).Value = 0.5499781580379022
$ws.Cells.Item(22, 12).Value = 0.2514182764704032
$ws.Cells.Item(22, 14).Value = 1.535830262570833
$ws.Cells.Item(22, 15).Value = 3.027562782080736

$ws.Cells.Item(23, 2).Value = 0.7602872114176478
$ws.Cells.Item(23, 3).Value = 0.1133182481183468
$ws.Cells.Item(23, 4).Value = 0.05268602524589738
$ws.Cells.Item(23, 5).Value = 0.1019421750801399
$ws.Cells.Item(23, 6).Value = 0.877842524539119
$ws.Cells.Item(23, 8).Value = 0.07973214163530429
$ws.Cells.Item(23, 9).Value = 0.7766568894895869
$ws.Cells.Item(23, 11).Value = 0.5317243329693326
$ws.Cells.Item(23, 12).Value = 0.2474291740567054
$ws.Cells.Item(23, 14).Value = 1.542063685384129
$ws.Cells.Item(23, 15).Value = 3.028206648975782

$ws.Cells.Item(24, 2).Value = 0.6826930923923555
$ws.Cells.Item(24, 3).Value = 0.1105160683478132
$ws.Cells.Item(24, 4).Value = 0.04875520140799239
$ws.Cells.Item(24, 5).Value = 0.1003963556700462
$ws.Cells.Item(24, 6).Value = 0.8720564830804065
$ws.Cells.Item(24, 8).Value = 0.07973214163530429
$ws.Cells.Item(24, 9).Value = 0.7820889833134856
$ws.Cells.Item(24, 11).Value = 0.4625834985604911
$ws.Cells.Item(24, 12).Value = 0.2325102803402643
$ws.Cells.Item(24, 14).Value = 1.566584440686333
$ws.Cells.Item(24, 15).Value = 3.033397773225772

$ws.Cells.Item(25, 2).Value = 0.5994920678619167
$ws.Cells.Item(25, 3).Value = 0.107451436002151
$ws.Cells.Item(25, 4).Value = 0.04445946247204091
$ws.Cells.Item(25, 5).Value = 0.09891083475780604
$ws.Cells.Item(25, 6).Value = 0.8676175942544475
$ws.Cells.Item(25, 8).Value = 0.07973214163530429
$ws.Cells.Item(25, 9).Value = 0.7893507933879746
$ws.Cells.Item(25, 11).Value = 0.3880098704330521
$ws.Cells.Item(25, 12).Value = 0.2168255955494942
$ws.Cells.Item(25, 14).Value = 1.594982092841631
$ws.Cells.Item(25, 15).Value = 3.044880393637783
